$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue $ws 'D2' '66.334.77'
Set-TextValue $ws 'E2' '  +0.39%  '
Set-TextValue $ws 'D3' '3.262.93'
Set-TextValue $ws 'E3' '  +2.83%  '
Set-TextValue $ws 'D4' '0.999'
Set-TextValue $ws 'E4' '  -0.05%  '
Set-TextValue $ws 'D5' '613.78'
Set-TextValue $ws 'E5' '  +0.91%  '
Set-TextValue $ws 'D6' '157.75'
Set-TextValue $ws 'E6' '  +2.06%  '
Set-TextValue $ws 'E7' '  +0.06%  '
Set-TextValue $ws 'D8' '3.262.65'
Set-TextValue $ws 'E8' '  +2.84%  '
Set-TextValue $ws 'D9' '0.545'
Set-TextValue $ws 'E9' '  +0.09%  '
Set-TextValue $ws 'E10' '  +1.73%  '
Set-TextValue $ws 'D11' '5.80'
Set-TextValue $ws 'E11' '  +2.07%  '
Set-TextValue $ws 'D12' '0.496'
Set-TextValue $ws 'E12' '  -4.24%  '
Set-TextValue $ws 'D13' '0.0000271'
Set-TextValue $ws 'E13' '  +1.02%  '
Set-TextValue $ws 'D14' '39.14'
Set-TextValue $ws 'E14' '  +1.87%  '
Set-TextValue $ws 'D15' '3.796.16'
Set-TextValue $ws 'E15' '  +2.78%  '
Set-TextValue $ws 'D16' '66.408.08'
Set-TextValue $ws 'E16' '  +0.41%  '
Set-TextValue $ws 'D17' '7.45'
Set-TextValue $ws 'E17' '  +0.28%  '
Set-TextValue $ws 'D18' '3.264.39'
Set-TextValue $ws 'E18' '  +3.12%  '
Set-TextValue $ws 'E19' '  +0.99%  '
Set-TextValue $ws 'D20' '505.37'
Set-TextValue $ws 'E20' '  -1.10%  '
Set-TextValue $ws 'D21' '15.47'
Set-TextValue $ws 'E21' '  +0.39%  '
Set-TextValue $ws 'D22' '0.756'
Set-TextValue $ws 'E22' '  +3.45%  '
Set-TextValue $ws 'D23' '8.11'
Set-TextValue $ws 'E23' '  +1.06%  '
Set-TextValue $ws 'D24' '14.67'
Set-TextValue $ws 'E24' '  -1.31%  '
Set-TextValue $ws 'D25' '87.40'
Set-TextValue $ws 'E25' '  +3.26%  '
Set-TextValue $ws 'E26' '  -0.06%  '
Set-TextValue $ws 'D27' '3.04'
Set-TextValue $ws 'E27' '  +0.86%  '
Set-TextValue $ws 'E28' '  +0.33%  '
Set-TextValue $ws 'D29' '2.39'
Set-TextValue $ws 'E29' '  +0.64%  '
Set-TextValue $ws 'D30' '0.131'
Set-TextValue $ws 'E30' '  +48.02%  '
Set-TextValue $ws 'D31' '7.10'
Set-TextValue $ws 'E31' '  -1.20%  '
Set-TextValue $ws 'E32' '  -4.27%  '
Set-TextValue $ws 'D33' '28.01'
Set-TextValue $ws 'E33' '  +0.06%  '
Set-TextValue $ws 'E34' '  -0.13%  '
Set-TextValue $ws 'E35' '  -3.75%  '
Set-TextValue $ws 'D36' '6.48'
Set-TextValue $ws 'E36' '  -0.62%  '
Set-TextValue $ws 'D37' '3.38'
Set-TextValue $ws 'E37' '  +19.28%  '
Set-TextValue $ws 'D38' '55.83'
Set-TextValue $ws 'E38' '  +1.45%  '
Set-TextValue $ws 'B39' 'PEPE'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws 'D39' '0.0₃0786'
Set-TextValue $ws 'E39' '  +15.08%  '
Set-TextValue $ws 'B40' 'Bittensor'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D40' '496.46'
Set-TextValue $ws 'E40' '  -1.22%  '
Set-TextValue $ws 'D41' '0.0422'
Set-TextValue $ws 'E41' '  +0.45%  '
Set-TextValue $ws 'E42' '  +1.03%  '
Set-TextValue $ws 'E43' '  +0.75%  '
Set-TextValue $ws 'B44' 'TheGraph'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws 'D44' '0.293'
Set-TextValue $ws 'E44' '  -1.26%  '
Set-TextValue $ws 'B45' 'Maker'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D45' '3.010.42'
Set-TextValue $ws 'E45' '  +6.67%  '
Set-TextValue $ws 'B46' 'Fetch.AI'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D46' '2.52'
Set-TextValue $ws 'E46' '  +3.56%  '
Set-TextValue $ws 'E47' '  +2.94%  '
Set-TextValue $ws 'E48' '  +5.44%  '
Set-TextValue $ws 'E49' '  +2.22%  '
Set-TextValue $ws 'E51' '  -3.47%  '
